# Updated queries for C3DC first half testcases.
#
# The SQL stored in the sheet's query cells joined tables using the bare
# "id" column (std.id / prt.id) aliased against "study.id" / "participant.id"
# quoted column names. Update every occurrence across the sheet to use the
# fully-qualified "study_id" / "participant_id" column names instead, e.g.:
#   df_participant prt ON std.id = prt."study.id"
#     -> df_participant prt ON std.study_id = prt."study.study_id"
#   df_diagnoses dgn ON prt.id = dgn."participant.id"
#     -> df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"
# (same pattern for df_treatments, df_treatment_resp, df_survival,
#  df_reference_files)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldJoinBlock = @"
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
"@

$newJoinBlock = @"
LEFT JOIN 
    df_participant prt ON std.study_id = prt."study.study_id"
LEFT JOIN 
    df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"
LEFT JOIN 
    df_treatments trt ON prt.participant_id = trt."participant.participant_id"
LEFT JOIN 
    df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"
LEFT JOIN 
    df_survival srv ON prt.participant_id = srv."participant.participant_id"
LEFT JOIN 
    df_reference_files rfs ON std.study_id = rfs."study.study_id"
"@

# The workbook has one query cell per tab, each embedding the same LEFT
# JOIN chain: the overall stats query (C2) plus one tab-query per table
# (B2 StudiesTab, B3 ParticipantsTab, B4 DiagnosisTab, B5 TreatmentTab,
# B6 TreatmentRespTab, B7 SurvivalTab).
$cellsToUpdate = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cellsToUpdate) {
    $range = $ws.Range($addr)
    $current = $range.Value()
    if ($current -ne $null -and $current.Contains($oldJoinBlock)) {
        $updated = $current.Replace($oldJoinBlock, $newJoinBlock)
        $range.Value = $updated
    }
}
